$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numbers as plain text (e.g. "578.32"). Setting
# Range.Value to a numeric-looking string auto-converts it to a real number,
# so for the price cells whose new value still parses as a plain number we
# briefly force a Text format on that single cell, assign the value, then
# ClearFormats() so the cell ends up back on the workbook default style
# (matching the source file, where none of these cells carry an explicit
# style). Non-numeric-looking text (URLs, coin names, multi-dot prices,
# percentages) round-trips as text natively and needs no special handling.

$ws.Range("D2").Value = "66.944.39"
$ws.Range("E2").Value = "  +4.17%  "
$ws.Range("D3").Value = "3.247.16"
$ws.Range("E3").Value = "  +2.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.91"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "3.244.56"
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.71"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.80%  "
$ws.Range("D13").Value = "3.810.53"
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.137"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.90"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").Value = "66.895.06"
$ws.Range("E16").Value = "  +4.11%  "
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").Value = "3.248.87"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("E19").Value = "  +1.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.31%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.73"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.508"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("D26").Value = "3.388.86"
$ws.Range("E26").Value = "  +2.18%  "
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.78"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.179"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "174.71"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.07%  "
$ws.Range("E36").Value = "  +2.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.75"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("E38").Value = "  +4.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.855"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.75%  "
$ws.Range("E40").Value = "  +9.11%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.76"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.49"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.81%  "
$ws.Range("D43").Value = "2.748.36"
$ws.Range("E43").Value = "  +2.88%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.56"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("E45").Value = "  +2.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.35"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("E47").Value = "  +3.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.65"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "335.58"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("E50").Value = "  +2.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.104"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.45%  "
